$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a blank row at the very top (old row1 "Filter Values" -> row2, etc.)
$ws.Rows("1:1").Insert()

# 2) Insert another blank row above the headers row (currently row3) so headers land on row4
$ws.Rows("3:3").Insert()

# 3) Insert a new column before column E so old E..H (now on rows 4-9) shift to F..I
$ws.Range("E1").EntireColumn.Insert()

# 4) New header row 3 (written first so shared-string indices land in the same
#    order as the source commit: Surname, Givenname, suburb, Pincode before
#    the "# of Records" family of strings)
$ws.Range("A3").Value = "Surname"
$ws.Range("B3").Value = "Givenname"
$ws.Range("C3").Value = "suburb"
$ws.Range("D3").Value = "Pincode"

# 5) New summary row 1
$ws.Range("A1").Value = "# of Records "
$ws.Range("B1").Value = 10000
$ws.Range("D1").Value = "# of Duplicates"
$ws.Range("E1").Value = 2362
$ws.Range("G1").Value = "# of unique records"
$ws.Range("H1").Value = 7638
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true
$ws.Range("H1").Font.Bold = $true

# 6) Fix row 9 (old row 7) values: 212 -> 112
$ws.Range("F9").Formula = "=1-(359+112)/10000"
$ws.Range("H9").Formula = "=1-(112/7638)"
$ws.Range("I9").Formula = "=1-(112/2115)"

# 7) New row 10
$ws.Range("A10").Value = 80
$ws.Range("B10").Value = 80
$ws.Range("C10").Value = 50
$ws.Range("D10").Value = 50
$ws.Range("F10").Formula = "=1-(212+198)/10000"
$ws.Range("G10").Formula = "=1-(212/2362)"
$ws.Range("H10").Formula = "=1-(198/7638)"
$ws.Range("I10").Formula = "=1-(198/1928)"
$ws.Range("F10:I10").NumberFormat = "0%"

# 8) Column widths
$ws.Columns("A").ColumnWidth = 11.54296875
$ws.Columns("B").ColumnWidth = 5.81640625
$ws.Columns("G").ColumnWidth = 17.1796875

# 9) Page orientation -> portrait (forces a pageSetup/printerSettings part)
$ws.PageSetup.Orientation = 1
